# Supervision section reformatted: merge the 3 existing supervisee rows
# down to 2 (Milena Vasquez-Amezquita + Francisco Javier Flores), drop the
# Julia Sanz-Vidania / Stirling / S Craig Roberts entry, expand the thesis
# descriptions, move the "Supervision conjunta con ..." notes to their own
# row beneath each supervisee, and trim the trailing blank rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Milena Vasquez-Amezquita (PhD in Neuroscience) ---
$ws.Range("A2").Value = "PhD in Neuroscience"
$ws.Range("B2").Value = "2015 - 2018"
$ws.Range("C2").Value = "\href{https://www.researchgate.net/profile/Milena-Vasquez-Amezquita}{Milena Vásquez-Amézquita}"
$ws.Range("D2").Value = "\href{https://www.uv.es/}{Universitat de València}, España"
$ws.Range("E2").Value = "Tésis \textbf{\textit{(Summa Cum Laude)}}: \textit{\href{http://hdl.handle.net/10550/67639}{Preferencias sexuales típicas y atípicas según sexo y edad de los estímulos: Utilidad de la técnica de rastreo ocular} [Typical and atypical sexual preferences according to sex and age of the stimuli: Usefulness of the eye tracking technique]}"
# Row 2 was already ht=75 and stays that way - leave the row height alone.

# --- Row 3: joint-supervision note for Milena, other cells blank ---
$ws.Range("A3").Value = $null
$ws.Range("B3").Value = $null
$ws.Range("C3").Value = $null
$ws.Range("D3").Value = $null
$ws.Range("E3").Value = "Supervisión conjunta con  Alicia Salvador"
# Row 3 used to be the 60pt Francisco row - drop back to the default height.
$ws.Rows(3).AutoFit()

# --- Row 4: Francisco Javier Flores (Professional Doctorate) ---
$ws.Range("A4").Value = "Professional Doctorate in Counselling Psychology"
$ws.Range("B4").Value = "2015 - 2018"
$ws.Range("C4").Value = "\href{https://www.researchgate.net/profile/Francisco-Flores-14}{Francisco Javier Flores}"
$ws.Range("D4").Value = "\href{https://www.uel.ac.uk/}{U. of East London}, Reino Unido"
$ws.Range("E4").Value = "Tésis: \textit{ What sense do people make of the functions of their ’behaviours that may be causing problems in their everyday life’? A hybrid deductive/inductive template analysis}"
$ws.Range("A4:E4").RowHeight = 60

# --- Row 5: joint-supervision note for Francisco, other cells blank ---
$ws.Range("A5").Value = $null
$ws.Range("B5").Value = $null
$ws.Range("C5").Value = $null
$ws.Range("D5").Value = $null
$ws.Range("E5").Value = "Supervisión conjunta con Lisa Chiara Fellin"
# Row 5 was already default height - make sure it stays that way.
$ws.Rows(5).AutoFit()

# E8 loses its bold-font override, back to the plain wrap style.
$ws.Range("E8").Font.Bold = $false

# Drop the now-unused trailing blank rows (16-21).
$ws.Range("A16:E21").EntireRow.Delete()

# Selection follows the data down to the new last row.
$ws.Range("A15").Select()
